# Apply cryptos list update (GitHub Actions style data refresh).
# Every changed cell is written with a leading apostrophe to force
# text storage (matches the source inlineStr cells, many of which look
# numeric, e.g. "1.000" / "25.089.41"), then the style is reset back to
# "Normal" so the quote-prefix flag does not leave a stray cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.089.41"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -2.99%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.652.54"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -4.99%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.12%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'236.93"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.83%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.17%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4801"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -8.15%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2621"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -4.59%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.05981"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -2.99%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07087"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.27%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'1.660.18"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -4.43%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'14.42"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -3.93%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.6192"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -3.73%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'4.588"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.62%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'72.93"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -5.77%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.09%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.17%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'25.080.53"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -3.10%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'11.34"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.37%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.000006525"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -3.71%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'4.410"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +3.24%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'1.870.75"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -4.65%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'8.449"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -2.24%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'5.256"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.22%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'133.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -3.91%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'14.75"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -2.90%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'1.389"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -8.52%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'1.703"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -3.71%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'101.49"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -3.52%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'3.816"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -3.49%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'0.07891"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -4.71%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'3.514"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -4.59%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.04597"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.57%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'2.604"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -1.42%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.9430"
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = "'0.5846"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -5.30%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'2.616"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -2.50%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -4.12%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = "'TrustWalletToken"
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'0.8392"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +12.87%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = "'PaxDollar"
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = "'1.001"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.21%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'1.831"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -5.27%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'98.29"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.69%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.3698"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -4.08%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'4.831"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -3.20%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.1126"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.24%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'6.042"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -3.32%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.05146"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -1.88%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'52.08"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -4.77%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'29.26"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -4.04%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.16%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'USDD"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'0.9976"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.03%  "
$ws.Range('E51').Style = 'Normal'
